$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Il34"
$ws.Cells.Item(2, 3).Value = "Ptprz1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.877313
$ws.Cells.Item(2, 8).Value = 5.631939
$ws.Cells.Item(2, 9).Value = 0.09253425024289272
$ws.Cells.Item(2, 10).Value = 0.0925342502428927
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.03995766666666666
$ws.Cells.Item(2, 14).Value = 0.119873
$ws.Cells.Item(2, 15).Value = 0.005314930928687666
$ws.Cells.Item(2, 16).Value = 0.005314930928687667
$ws.Cells.Item(2, 17).Value = 0.075013047083
$ws.Cells.Item(2, 18).Value = 0.675117423747
$ws.Cells.Item(2, 19).Value = 0.0004918131485788747
$ws.Cells.Item(2, 20).Value = 0.0004918131485788747

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Il34"
$ws.Cells.Item(3, 3).Value = "Ptprz1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.877313
$ws.Cells.Item(3, 8).Value = 5.631939
$ws.Cells.Item(3, 9).Value = 0.09253425024289272
$ws.Cells.Item(3, 10).Value = 0.0925342502428927
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.03069133333333333
$ws.Cells.Item(3, 14).Value = 0.092074
$ws.Cells.Item(3, 15).Value = 0.004082378436578614
$ws.Cells.Item(3, 16).Value = 0.004082378436578615
$ws.Cells.Item(3, 17).Value = 0.057617239054
$ws.Cells.Item(3, 18).Value = 0.5185551514860001
$ws.Cells.Item(3, 19).Value = 0.0003777598278365546
$ws.Cells.Item(3, 20).Value = 0.0003777598278365546

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Il34"
$ws.Cells.Item(4, 3).Value = "Ptprz1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.877313
$ws.Cells.Item(4, 8).Value = 5.631939
$ws.Cells.Item(4, 9).Value = 0.09253425024289272
$ws.Cells.Item(4, 10).Value = 0.0925342502428927
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 7.447354000000001
$ws.Cells.Item(4, 14).Value = 22.342062
$ws.Cells.Item(4, 15).Value = 0.9906026906347337
$ws.Cells.Item(4, 16).Value = 0.9906026906347338
$ws.Cells.Item(4, 17).Value = 13.981014479802
$ws.Cells.Item(4, 18).Value = 125.829130318218
$ws.Cells.Item(4, 19).Value = 0.09166467726647728
$ws.Cells.Item(4, 20).Value = 0.09166467726647728

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Il34"
$ws.Cells.Item(5, 3).Value = "Ptprz1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 3.680834
$ws.Cells.Item(5, 8).Value = 11.042502
$ws.Cells.Item(5, 9).Value = 0.1814312341407894
$ws.Cells.Item(5, 10).Value = 0.1814312341407894
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.03995766666666666
$ws.Cells.Item(5, 14).Value = 0.119873
$ws.Cells.Item(5, 15).Value = 0.005314930928687666
$ws.Cells.Item(5, 16).Value = 0.005314930928687667
$ws.Cells.Item(5, 17).Value = 0.1470775380273333
$ws.Cells.Item(5, 18).Value = 1.323697842246
$ws.Cells.Item(5, 19).Value = 0.0009642944777648551
$ws.Cells.Item(5, 20).Value = 0.0009642944777648553

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Il34"
$ws.Cells.Item(6, 3).Value = "Ptprz1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 3.680834
$ws.Cells.Item(6, 8).Value = 11.042502
$ws.Cells.Item(6, 9).Value = 0.1814312341407894
$ws.Cells.Item(6, 10).Value = 0.1814312341407894
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.03069133333333333
$ws.Cells.Item(6, 14).Value = 0.092074
$ws.Cells.Item(6, 15).Value = 0.004082378436578614
$ws.Cells.Item(6, 16).Value = 0.004082378436578615
$ws.Cells.Item(6, 17).Value = 0.1129697032386667
$ws.Cells.Item(6, 18).Value = 1.016727329148
$ws.Cells.Item(6, 19).Value = 0.0007406709579782043
$ws.Cells.Item(6, 20).Value = 0.0007406709579782044

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Il34"
$ws.Cells.Item(7, 3).Value = "Ptprz1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 3.680834
$ws.Cells.Item(7, 8).Value = 11.042502
$ws.Cells.Item(7, 9).Value = 0.1814312341407894
$ws.Cells.Item(7, 10).Value = 0.1814312341407894
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 7.447354000000001
$ws.Cells.Item(7, 14).Value = 22.342062
$ws.Cells.Item(7, 15).Value = 0.9906026906347337
$ws.Cells.Item(7, 16).Value = 0.9906026906347338
$ws.Cells.Item(7, 17).Value = 27.41247381323601
$ws.Cells.Item(7, 18).Value = 246.712264319124
$ws.Cells.Item(7, 19).Value = 0.1797262687050463
$ws.Cells.Item(7, 20).Value = 0.1797262687050464

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Il34"
$ws.Cells.Item(8, 3).Value = "Ptprz1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 14.72961666666667
$ws.Cells.Item(8, 8).Value = 44.18885
$ws.Cells.Item(8, 9).Value = 0.7260345156163179
$ws.Cells.Item(8, 10).Value = 0.7260345156163179
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.03995766666666666
$ws.Cells.Item(8, 14).Value = 0.119873
$ws.Cells.Item(8, 15).Value = 0.005314930928687666
$ws.Cells.Item(8, 16).Value = 0.005314930928687667
$ws.Cells.Item(8, 17).Value = 0.5885611128944443
$ws.Cells.Item(8, 18).Value = 5.29705001605
$ws.Cells.Item(8, 19).Value = 0.003858823302343936
$ws.Cells.Item(8, 20).Value = 0.003858823302343937

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Il34"
$ws.Cells.Item(9, 3).Value = "Ptprz1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 14.72961666666667
$ws.Cells.Item(9, 8).Value = 44.18885
$ws.Cells.Item(9, 9).Value = 0.7260345156163179
$ws.Cells.Item(9, 10).Value = 0.7260345156163179
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.03069133333333333
$ws.Cells.Item(9, 14).Value = 0.092074
$ws.Cells.Item(9, 15).Value = 0.004082378436578614
$ws.Cells.Item(9, 16).Value = 0.004082378436578615
$ws.Cells.Item(9, 17).Value = 0.4520715749888889
$ws.Cells.Item(9, 18).Value = 4.0686441749
$ws.Cells.Item(9, 19).Value = 0.002963947650763855
$ws.Cells.Item(9, 20).Value = 0.002963947650763856

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Il34"
$ws.Cells.Item(10, 3).Value = "Ptprz1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 14.72961666666667
$ws.Cells.Item(10, 8).Value = 44.18885
$ws.Cells.Item(10, 9).Value = 0.7260345156163179
$ws.Cells.Item(10, 10).Value = 0.7260345156163179
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 7.447354000000001
$ws.Cells.Item(10, 14).Value = 22.342062
$ws.Cells.Item(10, 15).Value = 0.9906026906347337
$ws.Cells.Item(10, 16).Value = 0.9906026906347338
$ws.Cells.Item(10, 17).Value = 109.6966696009667
$ws.Cells.Item(10, 18).Value = 987.2700264087001
$ws.Cells.Item(10, 19).Value = 0.7192117446632101
$ws.Cells.Item(10, 20).Value = 0.7192117446632103
